# Regenerate save_data: update column G (K, strike count) values for rows 2-37
# with the freshly recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 4
    4  = 4
    5  = 4
    6  = 3
    7  = 6
    8  = 3
    9  = 5
    10 = 5
    11 = 3
    12 = 3
    13 = 0
    14 = 1
    15 = 6
    16 = 5
    17 = 9
    18 = 2
    19 = 3
    20 = 6
    21 = 3
    22 = 2
    23 = 1
    24 = 3
    25 = 2
    26 = 0
    27 = 2
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
